$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.240.14"
$ws.Range("E2").Value = "  -5.22%  "
$ws.Range("D3").Value = "1.559.58"
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.19"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.38"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3408"
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.161"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07630"
$ws.Range("E11").Value = "  -5.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.36"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.019"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.917"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("D16").Value = "1.564.33"
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("E17").Value = "  -7.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.62"
$ws.Range("E18").Value = "  -5.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06723"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.225"
$ws.Range("E21").Value = "  -6.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.51"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5272"
$ws.Range("E23").Value = "  -7.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.93"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").Value = "22.251.29"
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.403"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.792"
$ws.Range("E27").Value = "  -6.66%  "
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.71"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.980"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.50"
$ws.Range("E31").Value = "  -4.64%  "
$ws.Range("D32").Value = "1.733.99"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.158"
$ws.Range("E33").Value = "  -10.31%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.015"
$ws.Range("E34").Value = "  -5.77%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.02"
$ws.Range("E36").Value = "  -10.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08496"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02535"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2310"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.473"
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.306"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.06380"
$ws.Range("E42").Value = "  -5.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.68"
$ws.Range("E43").Value = "  -9.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6326"
$ws.Range("E44").Value = "  -8.11%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.12"
$ws.Range("E45").Value = "  -8.86%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5967"
$ws.Range("E47").Value = "  -6.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.750"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.082"
$ws.Range("E49").Value = "  -7.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.262"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.33"
$ws.Range("E51").Value = "  -2.45%  "
